$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.404.02"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "2.363.31"
$ws.Range("E3").Value = "  +4.68%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("E5").Value = "  +1.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "232.23"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.74"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.79%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0946"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.18%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.94"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.42"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "2.704.60"
$ws.Range("E13").Value = "  +4.13%  "
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.45"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.27"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "2.354.63"
$ws.Range("E18").Value = "  +4.02%  "
$ws.Range("D19").Value = "43.332.41"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.58"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.24"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.46"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.93"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +20.00%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  +0.53%  "
$ws.Range("E27").Value = "  +0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.97"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.34"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.09"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.50"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.12%  "
$ws.Range("E32").Value = "  -7.26%  "
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.98"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0696"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.05"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.52"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +9.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.48"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.63"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("E40").Value = "  -2.73%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.97"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.08%  "
$ws.Range("B42").Value = "BinanceUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.80"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.94"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0948"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.36"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "1.441.75"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("B50").Value = "TerraClassic"
$ws.Range("C50").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000204"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -9.21%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.577.68"
$ws.Range("E51").Value = "  +4.30%  "
